# Re-applies the latest cryptos.xlsx price/volume refresh onto Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.502.30"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "3.492.41"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.213"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.77%  "
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("E12").Value = "  -4.12%  "
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "4.059.57"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "599.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("D16").Value = "69.623.68"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").Value = "3.493.65"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.80%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.06%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.628.32"
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.80%  "
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "507.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.29%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("D42").Value = "0.0₃0772"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0462"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("E47").Value = "  -4.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.07%  "
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("E51").Value = "  -7.99%  "
